$wb = $excel.ActiveWorkbook

# "Spain" is the last sheet and serves as the template for the new country
# tab. Clone it, rename the clone to "Turkey" and fill in the country
# specific values, mirroring how the other country tabs were produced.
$template = $wb.Worksheets.Item("Spain")
$template.Copy($null, $template)

$turkey = $wb.Worksheets.Item($template.Index + 1)
$turkey.Name = "Turkey"

$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3297"

# Column D / rows 3-4 re-wrap slightly differently than on the "Spain" tab
# once the country-specific strings are in place, so Excel widens column D
# and grows rows 3-4 to fit the wrapped "Constants"/"Input Value" labels.
$turkey.Columns.Item(4).ColumnWidth = 17.42
$turkey.Rows.Item(3).RowHeight = 28.8
$turkey.Rows.Item(4).RowHeight = 28.8

# Restore the previously active sheet's selection (it's no longer the
# active tab once "Turkey" is appended) and make "Turkey" the active sheet.
[void]$template.Range("A1:D10").Select()
$turkey.Activate()
[void]$turkey.Range("G12").Select()
